$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The document's lone "_GoBack" bookmark currently sits at the very
# end of the "...message originated from it." paragraph. The edit
# moves it: a new bold run "z" is added to the empty, bold-formatted
# paragraph that immediately precedes the "Setup" heading, and the
# "_GoBack" bookmark is relocated to sit right after that new run.
# The old bookmark location loses its bookmark tags entirely.
# ------------------------------------------------------------------

# Step 1: drop the bookmark from its current location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 2: locate the target paragraph - the empty bold paragraph that
# sits immediately before the "Setup" heading paragraph.
$rng = $d.Content
$rng.Find.Execute("Setup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingPara = $rng.Paragraphs(1)
$targetPara = $headingPara.Previous()
$target = $targetPara.Range
$target.Collapse(1)

# Step 3: type the new bold "z" run. A temporary trailing placeholder
# character ("X") is inserted too, so the bookmark can be anchored
# strictly between "z" and the placeholder - i.e. NOT sitting exactly
# on the paragraph-end boundary, which otherwise gets auto-expanded to
# span the whole paragraph when a bookmark is dropped right on it.
$target.InsertAfter("zX")
$target.Font.Bold = $true

$afterZ = $target.Duplicate
$afterZ.Collapse(1)
$afterZ.MoveStart(1, 1) | Out-Null
$afterZ.Collapse(0)
$d.Bookmarks.Add("_GoBack", $afterZ)

# Step 4: remove the temporary placeholder character, leaving the
# bookmark collapsed immediately after "z".
$placeholder = $afterZ.Duplicate
$placeholder.MoveEnd(1, 1) | Out-Null
$placeholder.Delete()
